$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new daily column "27-dec" before column EX ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Shift columns EX:GB (and everything in between) one column to the right,
# freeing up column EX for the newly observed "27-dec" day.
$wsPrix.Range("EX1").EntireColumn.Insert()

# Header cell for the new column (same header style as its neighbours).
$wsPrix.Range("EX1").Value = "27-dec"

# No data collected yet for this new day -> placeholder "-" like the other
# not-yet-populated columns.
$wsPrix.Range("EX2:EX25").Value = "-"

# --- Sheet "Gaz": append the newest daily price observation ---
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date column to stay plain text (matching every other row in
# column A) instead of letting Excel auto-convert it to a date serial.
$wsGaz.Range("A182").NumberFormat = "@"
$wsGaz.Range("A182").Value = "2025-12-25"
$wsGaz.Range("A182").ClearFormats()

$wsGaz.Range("B182").Value = 27.5
